# Roboflow Image Annotation 7/2/2025 - Good Night
# Fill in the weekly tracking row (row 49) of the "Avances Etiquetado
# Roboflow" table with the latest annotation progress figures, and
# update the sheet's selection to reflect where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49: Fecha, Imagenes sin etiquetar, Imagenes etiquetadas sin revisar,
# Imagenes rechazadas, Imagenes etiquetadas y revisadas faltando de subir,
# Imagenes etiquetadas revisadas y subidas, Notas
$ws.Range("D49").Value = 45695
$ws.Range("E49").Value = 192
$ws.Range("F49").Value = 734
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 1012
$ws.Range("J49").Value = "N/A"

# Move the selection / scroll position to match where the author ended up
$ws.Range("F51").Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
